$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 39

# Column A reuses the existing shared string used by other BMW listings
$ws.Cells.Item($row, 1).Value = "bmw e46 - od:50000, do:100000"

# New data for the newly scraped BMW 325tds listing
$ws.Cells.Item($row, 2).Value = "Bmw 325tds"
$ws.Cells.Item($row, 3).Value = "70 000 Kč"
# Column D reuses an existing "- [date]" shared string
$ws.Cells.Item($row, 4).Value = "- [8.2. 2024]"
$ws.Cells.Item($row, 5).Value = "02.12 2024 18:09"
$ws.Cells.Item($row, 6).Value = "Chrudim 537 01"
$ws.Cells.Item($row, 7).Value = "Prodám bmw e36 touring 2,5tds. Najeto 260 xxx Kastle bez rzi. Neschnilá jde zvedat za všechny zvedací body. Samozřejmě kosmetika je. Stk jsem dělal před zimou. Takže na necelé 2 roky klid. Bílá kůže Palubní počítač Všechny okna v elektrice. Tažné zařízení Na autě jsem dělal repas vstřikova ..."
$ws.Cells.Item($row, 8).Value = "https://www.bazos.cz/img/1t/624/181027624.jpg?t=1707408320"
$ws.Cells.Item($row, 9).Value = "https://auto.bazos.cz/inzerat/181027624/bmw-325tds.php"

# Columns B through I inherit the column-level style (index 1) when a new
# value is stored in a previously empty cell. The source row did not carry
# any explicit cell style, so copy the (unstyled) formatting from a cell
# that has no direct formatting (J1) onto the new row's cells to strip it
# back out again.
$ws.Cells.Item(1, 10).Copy() | Out-Null
$target = $ws.Range($ws.Cells.Item($row, 2), $ws.Cells.Item($row, 9))
$target.PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
